$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values (B2:E2)
$ws.Range("B2").Value = 10.751229216652829
$ws.Range("C2").Value = 9.0122213228944847
$ws.Range("D2").Value = 8.3612921469643791
$ws.Range("E2").Value = -0.12996849305077129

# Row 3 data values (B3:E3) - C3 is cleared, D3 is newly added
$ws.Range("B3").Value = 28.400687749377862
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 1.4618867507713229
$ws.Range("E3").Value = 2.1252554919599076

# Update the selection to match the new range
$ws.Range("B1:E3").Select() | Out-Null
